$wb = $excel.ActiveWorkbook

# Names (sheet1/2/3 use the full code+suffix in column A; sheet4 uses the
# bare row-code in column A, without the trailing lowercase suffix letter).
$namesFull = @("YYFWBz", "YYBNWf", "YYNWFw", "YYWFWb")
$namesCode = @("YYFWB", "YYNWF", "YYWFW")

function Remove-MatchingRows($ws, $names) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = $lastRow; $r -ge 1; $r--) {
        $v = $ws.Cells.Item($r, 1).Value()
        if ($names -contains $v) {
            $ws.Rows.Item($r).Delete()
        }
    }
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "rotation con1 set") {
        Remove-MatchingRows $ws $namesCode
    } else {
        Remove-MatchingRows $ws $namesFull
    }
}
